# Update "想去人数" (wanted-to-go count) figures across the four sheets to
# reflect the regenerated data (gh-pages output at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value  = 3345
$ws.Range("F5").Value  = 345
$ws.Range("F6").Value  = 7827
$ws.Range("F16").Value = 126
$ws.Range("F23").Value = 448
$ws.Range("F24").Value = 180
$ws.Range("F30").Value = 1115
$ws.Range("F37").Value = 1017
$ws.Range("F41").Value = 3372

# Sheet "演出"
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F24").Value = 76
$ws.Range("F25").Value = 6790

# Sheet "本地生活"
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value  = 2032
$ws.Range("F5").Value  = 1362
$ws.Range("F8").Value  = 2194
$ws.Range("F9").Value  = 9015
$ws.Range("F10").Value = 1210

# Sheet "全部类型"
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value  = 3345
$ws.Range("F3").Value  = 2032
$ws.Range("F4").Value  = 7827
$ws.Range("F5").Value  = 1362
$ws.Range("F7").Value  = 2194
$ws.Range("F9").Value  = 1210
$ws.Range("F16").Value = 126
$ws.Range("F23").Value = 180
$ws.Range("F29").Value = 1115
$ws.Range("F34").Value = 1017
$ws.Range("F38").Value = 3372
